# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Simple view-count (F column) increments ---
# Sheet 1
$ws1.Cells.Item(4, 6).Value = 3203
$ws1.Cells.Item(7, 6).Value = 323
$ws1.Cells.Item(8, 6).Value = 7492
$ws1.Cells.Item(11, 6).Value = 1218
$ws1.Cells.Item(13, 6).Value = 134
$ws1.Cells.Item(14, 6).Value = 608
$ws1.Cells.Item(15, 6).Value = 1058
$ws1.Cells.Item(16, 6).Value = 151
$ws1.Cells.Item(27, 6).Value = 1168
$ws1.Cells.Item(31, 6).Value = 402
$ws1.Cells.Item(33, 6).Value = 178
$ws1.Cells.Item(35, 6).Value = 336
$ws1.Cells.Item(36, 6).Value = 266
$ws1.Cells.Item(37, 6).Value = 950
$ws1.Cells.Item(38, 6).Value = 452
$ws1.Cells.Item(39, 6).Value = 75
$ws1.Cells.Item(40, 6).Value = 32
$ws1.Cells.Item(41, 6).Value = 249

# Sheet 2
$ws2.Cells.Item(6, 6).Value = 102
$ws2.Cells.Item(8, 6).Value = 364
$ws2.Cells.Item(9, 6).Value = 579
$ws2.Cells.Item(11, 6).Value = 62
$ws2.Cells.Item(14, 6).Value = 131
$ws2.Cells.Item(25, 6).Value = 3498
$ws2.Cells.Item(26, 6).Value = 3498
$ws2.Cells.Item(28, 6).Value = 44

# Sheet 3
$ws3.Cells.Item(6, 6).Value = 1913
$ws3.Cells.Item(8, 6).Value = 2978
$ws3.Cells.Item(9, 6).Value = 1182
$ws3.Cells.Item(12, 6).Value = 497
$ws3.Cells.Item(13, 6).Value = 1948
$ws3.Cells.Item(14, 6).Value = 8537
$ws3.Cells.Item(15, 6).Value = 725

# Sheet 4
$ws4.Cells.Item(3, 6).Value = 3203
$ws4.Cells.Item(5, 6).Value = 1913
$ws4.Cells.Item(6, 6).Value = 323
$ws4.Cells.Item(7, 6).Value = 1182
$ws4.Cells.Item(10, 6).Value = 497
$ws4.Cells.Item(11, 6).Value = 1218
$ws4.Cells.Item(12, 6).Value = 102
$ws4.Cells.Item(13, 6).Value = 134
$ws4.Cells.Item(14, 6).Value = 725
$ws4.Cells.Item(15, 6).Value = 364
$ws4.Cells.Item(16, 6).Value = 579
$ws4.Cells.Item(17, 6).Value = 579
$ws4.Cells.Item(19, 6).Value = 608
$ws4.Cells.Item(20, 6).Value = 1058
$ws4.Cells.Item(21, 6).Value = 62
$ws4.Cells.Item(24, 6).Value = 151
$ws4.Cells.Item(25, 6).Value = 131
$ws4.Cells.Item(36, 6).Value = 178
$ws4.Cells.Item(40, 6).Value = 336
$ws4.Cells.Item(41, 6).Value = 266
$ws4.Cells.Item(44, 6).Value = 452
$ws4.Cells.Item(45, 6).Value = 75
$ws4.Cells.Item(46, 6).Value = 32
$ws4.Cells.Item(47, 6).Value = 3498
$ws4.Cells.Item(49, 6).Value = 44

# --- Sheet 1 (展览): rows 17-24 cyclic reshuffle (new listing pushed in, 宫村优子 rotates to bottom) ---
$ws1.Cells.Item(17, 3).Value = '【大会员提前抢】上海 洛天依歌行宇宙·无限遨游 沉浸式体验展'
$ws1.Cells.Item(17, 4).Value = '中山北路3300号 上海月星环球港'
$ws1.Cells.Item(17, 5).Value = '2024.09.15 10:00-10.31 22:00'
$ws1.Cells.Item(17, 6).Value = 1221
$ws1.Cells.Item(17, 7).Value = 98
$ws1.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91175'
$ws1.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/ei9COXS41724405861343.jpeg'

$ws1.Cells.Item(18, 3).Value = '上海·GH·第五人格同人ONLY 1.0'
$ws1.Cells.Item(18, 4).Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws1.Cells.Item(18, 5).Value = '2024.09.15 10:00-09.15 17:00'
$ws1.Cells.Item(18, 6).Value = 300
$ws1.Cells.Item(18, 7).Value = 68
$ws1.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90638'
$ws1.Cells.Item(18, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/plVDxJKi1723102207272.jpeg'

$ws1.Cells.Item(19, 3).Value = '上海·SCGE动漫游戏嘉年华'
$ws1.Cells.Item(19, 4).Value = '军工路1076号 纪希片场(秀场)'
$ws1.Cells.Item(19, 5).Value = '2024.09.15 10:00-09.16 17:00'
$ws1.Cells.Item(19, 6).Value = 5887
$ws1.Cells.Item(19, 7).Value = 70
$ws1.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89993'
$ws1.Cells.Item(19, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/aIJyQziE1723434354531.jpeg'

$ws1.Cells.Item(20, 3).Value = '上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华'
$ws1.Cells.Item(20, 4).Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws1.Cells.Item(20, 5).Value = '2024.09.15 11:00-09.16 16:00'
$ws1.Cells.Item(20, 6).Value = 2328
$ws1.Cells.Item(20, 7).Value = 65.8
$ws1.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90990'
$ws1.Cells.Item(20, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/DutuUgvA1724127081751.jpeg'

$ws1.Cells.Item(21, 3).Value = '上海·原神ONLY逐月节·原神&崩铁&崩三&绝区零·同人动漫嘉年华'
$ws1.Cells.Item(21, 4).Value = '杨树浦路198号(金茂北外滩)B1层 Terra Park北外滩'
$ws1.Cells.Item(21, 5).Value = '2024.09.15 09:30-09.17 22:30'
$ws1.Cells.Item(21, 6).Value = 4071
$ws1.Cells.Item(21, 7).Value = 78
$ws1.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89712'
$ws1.Cells.Item(21, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/e9g9lWiy1721904672057.jpeg'

$ws1.Cells.Item(22, 3).Value = '上海·城市动漫节2th'
$ws1.Cells.Item(22, 4).Value = '西藏南路1号 上海大世界'
$ws1.Cells.Item(22, 5).Value = '2024.09.15 10:00-09.16 18:00'
$ws1.Cells.Item(22, 6).Value = 2267
$ws1.Cells.Item(22, 7).Value = 68
$ws1.Cells.Item(22, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89186'
$ws1.Cells.Item(22, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/dyznHqyF1723780926438.jpeg'

$ws1.Cells.Item(23, 3).Value = '上海·夜蓝诗·恋与深空同人only'
$ws1.Cells.Item(23, 4).Value = '莫干山路50号 M50创意园'
$ws1.Cells.Item(23, 5).Value = '2024.09.15 11:00-09.15 21:00'
$ws1.Cells.Item(23, 6).Value = 227
$ws1.Cells.Item(23, 7).Value = 88
$ws1.Cells.Item(23, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90729'
$ws1.Cells.Item(23, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/dBiBf2Ac1723543844923.jpeg'

$ws1.Cells.Item(24, 3).Value = '上海·宫村优子粉丝见面会'
$ws1.Cells.Item(24, 4).Value = '西藏南路1号 上海大世界'
$ws1.Cells.Item(24, 5).Value = '2024.09.15 12:00-09.16 17:30'
$ws1.Cells.Item(24, 6).Value = 58
$ws1.Cells.Item(24, 7).Value = 198
$ws1.Cells.Item(24, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91139'
$ws1.Cells.Item(24, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/pk4s8Bxs1724644287023.jpeg'

# --- Sheet 4 (全部类型): rows 26-31 cyclic reshuffle (mirrors Sheet 1's block) ---
$ws4.Cells.Item(26, 3).Value = '上海·GH·第五人格同人ONLY 1.0'
$ws4.Cells.Item(26, 4).Value = '逸仙路301号靠纪念路路口 上海宝丰联大酒店'
$ws4.Cells.Item(26, 5).Value = '2024.09.15 10:00-09.15 17:00'
$ws4.Cells.Item(26, 6).Value = 300
$ws4.Cells.Item(26, 7).Value = 68
$ws4.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90638'
$ws4.Cells.Item(26, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/plVDxJKi1723102207272.jpeg'

$ws4.Cells.Item(27, 3).Value = '上海·SCGE动漫游戏嘉年华'
$ws4.Cells.Item(27, 4).Value = '军工路1076号 纪希片场(秀场)'
$ws4.Cells.Item(27, 5).Value = '2024.09.15 10:00-09.16 17:00'
$ws4.Cells.Item(27, 6).Value = 5887
$ws4.Cells.Item(27, 7).Value = 70
$ws4.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89993'
$ws4.Cells.Item(27, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/aIJyQziE1723434354531.jpeg'

$ws4.Cells.Item(28, 3).Value = '上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华'
$ws4.Cells.Item(28, 4).Value = '盈浦街道淀山浦社区淀山湖大道851号青浦万达茂F3 万达汽车乐园(青浦万达茂店)'
$ws4.Cells.Item(28, 5).Value = '2024.09.15 11:00-09.16 16:00'
$ws4.Cells.Item(28, 6).Value = 2328
$ws4.Cells.Item(28, 7).Value = 65.8
$ws4.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90990'
$ws4.Cells.Item(28, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/DutuUgvA1724127081751.jpeg'

$ws4.Cells.Item(29, 3).Value = '上海·原神ONLY逐月节·原神&崩铁&崩三&绝区零·同人动漫嘉年华'
$ws4.Cells.Item(29, 4).Value = '杨树浦路198号(金茂北外滩)B1层 Terra Park北外滩'
$ws4.Cells.Item(29, 5).Value = '2024.09.15 09:30-09.17 22:30'
$ws4.Cells.Item(29, 6).Value = 4071
$ws4.Cells.Item(29, 7).Value = 78
$ws4.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89712'
$ws4.Cells.Item(29, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/e9g9lWiy1721904672057.jpeg'

$ws4.Cells.Item(30, 3).Value = '上海·城市动漫节2th'
$ws4.Cells.Item(30, 4).Value = '西藏南路1号 上海大世界'
$ws4.Cells.Item(30, 5).Value = '2024.09.15 10:00-09.16 18:00'
$ws4.Cells.Item(30, 6).Value = 2267
$ws4.Cells.Item(30, 7).Value = 68
$ws4.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89186'
$ws4.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/dyznHqyF1723780926438.jpeg'

$ws4.Cells.Item(31, 3).Value = '上海·宫村优子粉丝见面会'
$ws4.Cells.Item(31, 4).Value = '西藏南路1号 上海大世界'
$ws4.Cells.Item(31, 5).Value = '2024.09.15 12:00-09.16 17:30'
$ws4.Cells.Item(31, 6).Value = 58
$ws4.Cells.Item(31, 7).Value = 198
$ws4.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=91139'
$ws4.Cells.Item(31, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/pk4s8Bxs1724644287023.jpeg'

Write-Output "edit applied"